$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.07901576910072249
$ws.Range("J2").Value = 0.07901576910072249
$ws.Range("M2").Value = 0.4468086666666666
$ws.Range("N2").Value = 1.340426
$ws.Range("O2").Value = 0.1280124579567358
$ws.Range("P2").Value = 0.1280124579567358
$ws.Range("Q2").Value = 0.0439586749251111
$ws.Range("R2").Value = 0.3956280743259999
$ws.Range("S2").Value = 0.01011500281992538
$ws.Range("T2").Value = 0.01011500281992538

$ws.Range("I3").Value = 0.07901576910072249
$ws.Range("J3").Value = 0.07901576910072249
$ws.Range("O3").Value = 0.1587947312683464
$ws.Range("P3").Value = 0.1587947312683464
$ws.Range("S3").Value = 0.01254728782031094
$ws.Range("T3").Value = 0.01254728782031094

$ws.Range("I4").Value = 0.07901576910072249
$ws.Range("J4").Value = 0.07901576910072249
$ws.Range("M4").Value = 1.022202333333333
$ws.Range("N4").Value = 3.066607
$ws.Range("O4").Value = 0.2928650292200627
$ws.Range("P4").Value = 0.2928650292200626
$ws.Range("Q4").Value = 0.1005680136285555
$ws.Range("R4").Value = 0.9051121226569999
$ws.Range("S4").Value = 0.02314095552652882
$ws.Range("T4").Value = 0.02314095552652881

$ws.Range("I5").Value = 0.07901576910072249
$ws.Range("J5").Value = 0.07901576910072249
$ws.Range("M5").Value = 0.3652233333333334
$ws.Range("N5").Value = 1.09567
$ws.Range("O5").Value = 0.1046379358572997
$ws.Range("P5").Value = 0.1046379358572996
$ws.Range("Q5").Value = 0.03593201068555556
$ws.Range("R5").Value = 0.32338809617
$ws.Range("S5").Value = 0.0082680469788766
$ws.Range("T5").Value = 0.008268046978876597

$ws.Range("I6").Value = 0.07901576910072249
$ws.Range("J6").Value = 0.07901576910072249
$ws.Range("M6").Value = 0.5234373333333334
$ws.Range("N6").Value = 1.570312
$ws.Range("O6").Value = 0.1499668753657104
$ws.Range("P6").Value = 0.1499668753657104
$ws.Range("Q6").Value = 0.05149768412355556
$ws.Range("R6").Value = 0.4634791571120001
$ws.Range("S6").Value = 0.0118497479966538
$ws.Range("T6").Value = 0.0118497479966538

$ws.Range("I7").Value = 0.07901576910072249
$ws.Range("J7").Value = 0.07901576910072249
$ws.Range("M7").Value = 0.5784316666666666
$ws.Range("N7").Value = 1.735295
$ws.Range("O7").Value = 0.1657229703318451
$ws.Range("P7").Value = 0.1657229703318451
$ws.Range("Q7").Value = 0.05690822828277776
$ws.Range("R7").Value = 0.512174054545
$ws.Range("S7").Value = 0.01309472795842696
$ws.Range("T7").Value = 0.01309472795842696

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1928733333333333
$ws.Range("H8").Value = 0.5786199999999999
$ws.Range("I8").Value = 0.1549041145619023
$ws.Range("J8").Value = 0.1549041145619023
$ws.Range("M8").Value = 0.4468086666666666
$ws.Range("N8").Value = 1.340426
$ws.Range("O8").Value = 0.1280124579567358
$ws.Range("P8").Value = 0.1280124579567358
$ws.Range("Q8").Value = 0.0861774769022222
$ws.Range("R8").Value = 0.7755972921199998
$ws.Range("S8").Value = 0.01982965645268091
$ws.Range("T8").Value = 0.01982965645268091

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1928733333333333
$ws.Range("H9").Value = 0.5786199999999999
$ws.Range("I9").Value = 0.1549041145619023
$ws.Range("J9").Value = 0.1549041145619023
$ws.Range("O9").Value = 0.1587947312683464
$ws.Range("P9").Value = 0.1587947312683464
$ws.Range("Q9").Value = 0.1068999807088889
$ws.Range("R9").Value = 0.9620998263799999
$ws.Range("S9").Value = 0.02459795724421843
$ws.Range("T9").Value = 0.02459795724421843

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1928733333333333
$ws.Range("H10").Value = 0.5786199999999999
$ws.Range("I10").Value = 0.1549041145619023
$ws.Range("J10").Value = 0.1549041145619023
$ws.Range("M10").Value = 1.022202333333333
$ws.Range("N10").Value = 3.066607
$ws.Range("O10").Value = 0.2928650292200627
$ws.Range("P10").Value = 0.2928650292200626
$ws.Range("Q10").Value = 0.1971555713711111
$ws.Range("R10").Value = 1.77440014234
$ws.Range("S10").Value = 0.04536599803747947
$ws.Range("T10").Value = 0.04536599803747946

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1928733333333333
$ws.Range("H11").Value = 0.5786199999999999
$ws.Range("I11").Value = 0.1549041145619023
$ws.Range("J11").Value = 0.1549041145619023
$ws.Range("M11").Value = 0.3652233333333334
$ws.Range("N11").Value = 1.09567
$ws.Range("O11").Value = 0.1046379358572997
$ws.Range("P11").Value = 0.1046379358572996
$ws.Range("Q11").Value = 0.07044184171111112
$ws.Range("R11").Value = 0.6339765754
$ws.Range("S11").Value = 0.01620884680356013
$ws.Range("T11").Value = 0.01620884680356013

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1928733333333333
$ws.Range("H12").Value = 0.5786199999999999
$ws.Range("I12").Value = 0.1549041145619023
$ws.Range("J12").Value = 0.1549041145619023
$ws.Range("M12").Value = 0.5234373333333334
$ws.Range("N12").Value = 1.570312
$ws.Range("O12").Value = 0.1499668753657104
$ws.Range("P12").Value = 0.1499668753657104
$ws.Range("Q12").Value = 0.1009571032711111
$ws.Range("R12").Value = 0.90861392944
$ws.Range("S12").Value = 0.02323048604214054
$ws.Range("T12").Value = 0.02323048604214054

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1928733333333333
$ws.Range("H13").Value = 0.5786199999999999
$ws.Range("I13").Value = 0.1549041145619023
$ws.Range("J13").Value = 0.1549041145619023
$ws.Range("M13").Value = 0.5784316666666666
$ws.Range("N13").Value = 1.735295
$ws.Range("O13").Value = 0.1657229703318451
$ws.Range("P13").Value = 0.1657229703318451
$ws.Range("Q13").Value = 0.1115640436555555
$ws.Range("R13").Value = 1.0040763929
$ws.Range("S13").Value = 0.02567116998182288
$ws.Range("T13").Value = 0.02567116998182288

$ws.Range("G14").Value = 0.9538573333333332
$ws.Range("H14").Value = 2.861572
$ws.Range("I14").Value = 0.7660801163373752
$ws.Range("J14").Value = 0.7660801163373753
$ws.Range("M14").Value = 0.4468086666666666
$ws.Range("N14").Value = 1.340426
$ws.Range("O14").Value = 0.1280124579567358
$ws.Range("P14").Value = 0.1280124579567358
$ws.Range("Q14").Value = 0.4261917232968888
$ws.Range("R14").Value = 3.835725509671999
$ws.Range("S14").Value = 0.09806779868412953
$ws.Range("T14").Value = 0.09806779868412951

$ws.Range("G15").Value = 0.9538573333333332
$ws.Range("H15").Value = 2.861572
$ws.Range("I15").Value = 0.7660801163373752
$ws.Range("J15").Value = 0.7660801163373753
$ws.Range("O15").Value = 0.1587947312683464
$ws.Range("P15").Value = 0.1587947312683464
$ws.Range("Q15").Value = 0.5286751090475554
$ws.Range("R15").Value = 4.758075981428
$ws.Range("S15").Value = 0.121649486203817
$ws.Range("T15").Value = 0.121649486203817

$ws.Range("G16").Value = 0.9538573333333332
$ws.Range("H16").Value = 2.861572
$ws.Range("I16").Value = 0.7660801163373752
$ws.Range("J16").Value = 0.7660801163373753
$ws.Range("M16").Value = 1.022202333333333
$ws.Range("N16").Value = 3.066607
$ws.Range("O16").Value = 0.2928650292200627
$ws.Range("P16").Value = 0.2928650292200626
$ws.Range("Q16").Value = 0.9750351918004442
$ws.Range("R16").Value = 8.775316726203998
$ws.Range("S16").Value = 0.2243580756560544
$ws.Range("T16").Value = 0.2243580756560544

$ws.Range("G17").Value = 0.9538573333333332
$ws.Range("H17").Value = 2.861572
$ws.Range("I17").Value = 0.7660801163373752
$ws.Range("J17").Value = 0.7660801163373753
$ws.Range("M17").Value = 0.3652233333333334
$ws.Range("N17").Value = 1.09567
$ws.Range("O17").Value = 0.1046379358572997
$ws.Range("P17").Value = 0.1046379358572996
$ws.Range("Q17").Value = 0.3483709548044445
$ws.Range("R17").Value = 3.13533859324
$ws.Range("S17").Value = 0.08016104207486292
$ws.Range("T17").Value = 0.08016104207486291

$ws.Range("G18").Value = 0.9538573333333332
$ws.Range("H18").Value = 2.861572
$ws.Range("I18").Value = 0.7660801163373752
$ws.Range("J18").Value = 0.7660801163373753
$ws.Range("M18").Value = 0.5234373333333334
$ws.Range("N18").Value = 1.570312
$ws.Range("O18").Value = 0.1499668753657104
$ws.Range("P18").Value = 0.1499668753657104
$ws.Range("Q18").Value = 0.4992845389404444
$ws.Range("R18").Value = 4.493560850464
$ws.Range("S18").Value = 0.1148866413269161
$ws.Range("T18").Value = 0.1148866413269161

$ws.Range("G19").Value = 0.9538573333333332
$ws.Range("H19").Value = 2.861572
$ws.Range("I19").Value = 0.7660801163373752
$ws.Range("J19").Value = 0.7660801163373753
$ws.Range("M19").Value = 0.5784316666666666
$ws.Range("N19").Value = 1.735295
$ws.Range("O19").Value = 0.1657229703318451
$ws.Range("P19").Value = 0.1657229703318451
$ws.Range("Q19").Value = 0.551741287082222
$ws.Range("R19").Value = 4.965671583739999
$ws.Range("S19").Value = 0.1269570723915953
$ws.Range("T19").Value = 0.1269570723915953
